# Testidokumentaatio.xlsx edit
# Adds four new integration-test rows (18-21, sheet rows 19-22) to the
# "Testaus" sheet's Table1, matching:
#   GlobalExceptionHandlerTest, AddAndRemovePermissionToRoleTest,
#   SearchSaleTest, MarkTicketAsUsedByBarcodeTest

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testaus")

# --- Expand Table1 by four rows (A1:J18 -> A1:J22) -------------------------
$lo = $ws.ListObjects.Item("Table1")
1..4 | ForEach-Object { $lo.ListRows.Add() | Out-Null }

# --- New row content ---------------------------------------------------
# columns: A=ID, B=Nimi, C=Tyyppi, D=Kohde, E=Toimenpiteet, F=Odotettu tulos,
#          G=Tekijä, H=Status, I=Kommentti, J=Testin tyypin valinnan perustelut
$rows = @(
    @{ Row = 19; Id = 18; Nimi = "GlobalExceptionHandlerTest";        Status = "OK" },
    @{ Row = 20; Id = 19; Nimi = "AddAndRemovePermissionToRoleTest";  Status = "OK" },
    @{ Row = 21; Id = 20; Nimi = "SearchSaleTest";                    Status = "IN PROGRESS" },
    @{ Row = 22; Id = 21; Nimi = "MarkTicketAsUsedByBarcodeTest";     Status = "IN PROGRESS" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Range("A$rowNum").Value = $r.Id
    $ws.Range("B$rowNum").Value = $r.Nimi
    $ws.Range("G$rowNum").Value = "Julia"
    $ws.Range("H$rowNum").Value = $r.Status

    # Apply the same cell styling pattern used by the sheet's other data
    # rows: A/B = vertical-top; C/D/F/G/H/I = left+top; E/J = left+top+wrap.
    $ws.Range("A$rowNum`:B$rowNum").VerticalAlignment = -4160

    $ws.Range("C$rowNum").VerticalAlignment = -4160
    $ws.Range("C$rowNum").HorizontalAlignment = -4131
    $ws.Range("D$rowNum").VerticalAlignment = -4160
    $ws.Range("D$rowNum").HorizontalAlignment = -4131
    $ws.Range("F$rowNum").VerticalAlignment = -4160
    $ws.Range("F$rowNum").HorizontalAlignment = -4131
    $ws.Range("G$rowNum").VerticalAlignment = -4160
    $ws.Range("G$rowNum").HorizontalAlignment = -4131
    $ws.Range("H$rowNum").VerticalAlignment = -4160
    $ws.Range("H$rowNum").HorizontalAlignment = -4131
    $ws.Range("I$rowNum").VerticalAlignment = -4160
    $ws.Range("I$rowNum").HorizontalAlignment = -4131

    $ws.Range("E$rowNum").VerticalAlignment = -4160
    $ws.Range("E$rowNum").HorizontalAlignment = -4131
    $ws.Range("E$rowNum").WrapText = $true
    $ws.Range("J$rowNum").VerticalAlignment = -4160
    $ws.Range("J$rowNum").HorizontalAlignment = -4131
    $ws.Range("J$rowNum").WrapText = $true
}

# --- Extend the Status column conditional formatting / data validation -----
# so the new rows (H19:H22) keep the same "OK/NOK/TODO/IN PROGRESS" list
# validation and colour-coding as the rest of the column.
$ws.Range("H2:H18").Copy() | Out-Null
$ws.Range("H19:H22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "Added 4 new Table1 rows (19-22) to Testaus sheet."
